$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Unhide columns B:D (RegionID, Cases, Recovered) so the refreshed temporal /
# mobility data (now extending through mid-September) is visible again. Their
# stored widths already reflect the correct best-fit sizing, so unhiding is
# all that's needed to bring them back on screen.
$ws.Range("B:D").EntireColumn.Hidden = $false

# Point the active selection at the newly-updated data block at the bottom of
# the sheet.
$ws.Range("A31:E35").Select()
